$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29
$ws.Range("N29").Value = 46031
$ws.Range("Q29").Value = 2.24
$ws.Range("R29").Value = 2.23
$ws.Range("S29").Value = 2.24
$ws.Range("T29").Value = 2.24
$ws.Range("U29").Value = 2.23

# Row 30
$ws.Range("N30").Value = 46031
$ws.Range("Q30").Value = 2.28
$ws.Range("R30").Value = 2.27
$ws.Range("S30").Value = 2.27
$ws.Range("T30").Value = 2.27
$ws.Range("U30").Value = 2.26

# Row 47
$ws.Range("N47").Value = 46030

# Row 48
$ws.Range("N48").Value = 46030
$ws.Range("Q48").Value = 3.49
$ws.Range("R48").Value = 3.47
$ws.Range("S48").Value = 3.47
$ws.Range("T48").Value = 3.46
$ws.Range("U48").Value = 3.47

# Row 49
$ws.Range("N49").Value = 46030
$ws.Range("Q49").Value = 3.74
$ws.Range("R49").Value = 3.7
$ws.Range("S49").Value = 3.72
$ws.Range("T49").Value = 3.71
$ws.Range("U49").Value = 3.74

# Row 50
$ws.Range("N50").Value = 46030
$ws.Range("Q50").Value = 4.19
$ws.Range("R50").Value = 4.15
$ws.Range("S50").Value = 4.18
$ws.Range("T50").Value = 4.17
$ws.Range("U50").Value = 4.19

# Row 52
$ws.Range("N52").Value = 46030
$ws.Range("Q52").Value = 5.92
$ws.Range("R52").Value = 5.88
$ws.Range("S52").Value = 5.92
$ws.Range("T52").Value = 5.92
$ws.Range("U52").Value = 5.93
